# Cat, Parrot, & a Bag of Seed - add the final step paragraph.
#
# We append a new ListParagraph (same numbering as the preceding list items)
# containing the final step text, with "river bank" wrapped in
# proofErr gramStart/gramEnd markers (mirroring Word's grammar-checker
# output), and we relocate the document's "_GoBack" bookmark so it still
# marks the point of the very last edit, i.e. the end of the new paragraph.

$d = $word.ActiveDocument

# Remove the existing _GoBack bookmark; we'll recreate it at the new
# insertion point once the new paragraph has been added.
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# Build the new paragraph's OOXML: same list style/numbering as the other
# steps, three runs with proofErr markers bracketing "river bank", and a
# couple of throwaway trailing characters ("~~") that we use purely as a
# safe landing spot for the new bookmark before trimming them back off
# (this runtime mishandles bookmarks placed exactly at the very end of the
# document, so we give it a little runway past the real content first).
$apos = [char]0x2019
$newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
  '<w:r><w:t xml:space="preserve">Since we know the cat won' + $apos + 't eat the seed, we can leave those two together and go get the parrot and bring it back across the </w:t></w:r>' +
  '<w:proofErr w:type="gramStart"/>' +
  '<w:r><w:t>river bank</w:t></w:r>' +
  '<w:proofErr w:type="gramEnd"/>' +
  '<w:r><w:t>.~~</w:t></w:r>' +
  '</w:p>'

$insertAt = $d.Content
$insertAt.Collapse(0)
[void]$insertAt.InsertXML($newParaXml)

# The new paragraph (with its "~~" placeholder tail) is now the last
# paragraph in the document. Find the position right after the final
# "." -- i.e. right before the two placeholder characters -- and anchor
# the recreated _GoBack bookmark there.
$docEnd = $d.Content.End
$bookmarkPos = $docEnd - 3
$bookmarkRange = $d.Range($bookmarkPos, $bookmarkPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

# Trim the two placeholder characters back off, one at a time, working
# from the end backwards so neither delete touches the bookmark's own
# position directly (which would otherwise corrupt/displace it).
$ce = $d.Content.End
$d.Range($ce - 2, $ce - 1).Delete()
$ce = $d.Content.End
$d.Range($ce - 2, $ce - 1).Delete()
